$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab/sheet name to "Session"
$ws.Name = "Session"

# Remove the second row (data row) entirely, shrinking the used range back to A1:F1
$ws.Rows.Item(2).Delete()
